# Update Sidebar Pasien and Pendaftaran
# - Sheet1: rename content to "Daftar Barang", set column A width
# - Add 3 new sheets: "Daftar Satuan Barang", "Daftar Kategori Barang", "Daftar Cabang"

$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "Daftar Barang"
$ws1.Columns.Item(1).ColumnWidth = 11.5

# --- Sheet2: Daftar Satuan Barang ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Daftar Satuan Barang"
$ws2.Range("A1").Value = "No"
$ws2.Range("B1").Value = "Satuan Barang"
$ws2.Columns.Item(1).ColumnWidth = 2.6666666666666665
$ws2.Columns.Item(2).ColumnWidth = 11.833333333333334
[void]$ws2.Range("C4").Select()

# --- Sheet3: Daftar Kategori Barang ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Daftar Kategori Barang"
$ws3.Range("A1").Value = "No"
$ws3.Range("B1").Value = "Kategori Barang"
$ws3.Columns.Item(1).ColumnWidth = 2.6666666666666665
$ws3.Columns.Item(2).ColumnWidth = 13.166666666666666
[void]$ws3.Range("C6").Select()

# --- Sheet4: Daftar Cabang ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Daftar Cabang"
$ws4.Range("A1").Value = "No"
$ws4.Range("B1").Value = "Daftar Cabang"
$ws4.Columns.Item(1).ColumnWidth = 2.6666666666666665
$ws4.Columns.Item(2).ColumnWidth = 11.833333333333334
[void]$ws4.Range("I13").Select()

# restore original active sheet/tab selection
$ws1.Activate()
